# Add slides for two sessions (09_policy and the new 10_hackathon session)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (09_policy): add Code and Data slide references
$ws.Range("F10").Value = "09_policy.R"
$ws.Range("G10").Value = "09_policy.RData"

# Row 13: "Round-up" session becomes "Hackathon" with its own slides
$ws.Range("D13").Value = "Hackathon"
$ws.Range("E13").Value = "10_hackathon"

# Update the active selection to match the author's final cursor position
$ws.Range("E14").Select()
